$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(4.262580982041177, 4.031621103296338, -5.649476117924912, 9.187166381605149)"
$ws.Range("C2").Value = "NIG(0.6715896618005597, 0.402211843272708, 6.491386786046982, 5.0720989498947)"
$ws.Range("D2").Value = "NCT(3.1810663867728217, 1.606193081331272, -1.136299839010911, 2.5059020270901167)"
$ws.Range("E2").Value = "NIG(1.153524607450748, 0.8802659054066564, 4.520701796837448, 5.551211290867952)"
